# This script reorders the "Recorded By" (column G) values for specific
# rows in the "Session Analysis Results" sheet, swapping the order of the
# comma-separated names/emails while keeping the same set of values.
#
#   "dnasr281@gmail.com, System"              -> "System, dnasr281@gmail.com"
#   "system, backup@backdoor.com, System"     -> "backup@backdoor.com, system, System"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# Rows whose "Recorded By" value is "dnasr281@gmail.com, System"
$rowsDnasr = @(10,11,12,13,14,15,17,18,19,20,21,22,24,26,36,37,38,39,40,41,43,44,45,46,47,48,50,52,62,63,64,65,66,67,69,70,71,72,73,74,76,78,83,84,85,86,90,92,93,94,96,99,101,109,110,111,112,116,118,119,120,122,125,127,135,136,137,138,142,144,145,146,148,151,153)

foreach ($r in $rowsDnasr) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -eq "dnasr281@gmail.com, System") {
        $cell.Value2 = "System, dnasr281@gmail.com"
    }
}

# Rows whose "Recorded By" value is "system, backup@backdoor.com, System"
$rowsSystem = @(2,28,54)

foreach ($r in $rowsSystem) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -eq "system, backup@backdoor.com, System") {
        $cell.Value2 = "backup@backdoor.com, system, System"
    }
}
